# Remove Marky Marc, and add flight details to reports
# (this workbook only covers the "remove Marky Marc" portion of that commit;
#  the roster table also drops the unrelated "Andre" solo row per the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# Find the worksheet rows for the two roster entries being removed by reading
# the table's Name column (column A of each list row), so this is resilient
# to the table not being in the exact position we expect.
$rowsToDelete = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $tbl.ListRows.Count; $i++) {
    $listRow = $tbl.ListRows.Item($i)
    $name = $listRow.Range.Cells.Item(1, 1).Value2
    if ($name -eq "Marky Marc" -or $name -eq "Andre") {
        [void]$rowsToDelete.Add($listRow.Range.Row)
    }
}

# Delete highest row number first so the remaining row numbers stay valid.
$sortedRows = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRows) {
    $ws.Rows($r).Delete()
}

# Mirror the author's final cursor position: the (now shifted) last data row.
$lastRow = $tbl.Range.Row + $tbl.ListRows.Count
[void]$ws.Range("A" + $lastRow + ":XFD" + $lastRow).Select()
